$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace header "Types" with "The closest predicted product" and drop "Angle" column header
$ws.Range("C1").Value = "The closest predicted product"
$ws.Range("D1").ClearContents()

# Fill in predicted-product SMILES values in column C
$ws.Range("C2").Value = "OCCCc1ccccc1"
$ws.Range("C3").Value = "OCCCc1cccc(c1)C(=O)C"
$ws.Range("C4").Value = "OCC(c1ccccc1CC)C"
$ws.Range("C5").Value = "OCC(c1ccc(cc1)CC)C"
$ws.Range("C6").Value = "OCCc1ccc(cc1)CC"

# Widen column C to fit the new, longer content
$ws.Columns.Item(3).ColumnWidth = 28.15

# Update the active selection to match the saved view state
$ws.Range("C10").Select()
